$wb = $excel.ActiveWorkbook

# --- Sheet: VENTAS POR GRUPO ---
$wsVentasGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsVentasGrupo.Range("D29").Value = 2260.22
$wsVentasGrupo.Range("K29").Value = 456.84
$wsVentasGrupo.Range("D53").Value = "6 de 51"
$wsVentasGrupo.Range("K53").Value = "2 de 51"

# --- Sheet: VENTA MENSUAL ---
$wsVentaMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsVentaMensual.Range("F29").Value = 2809.53
$wsVentaMensual.Range("F53").Value = 46225.1

# --- Sheet: CUMPLIMIENTO MENSUAL ---
$wsCumplimiento = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# Column F width: raw OOXML width 25 -> 24 (ColumnWidth property adds ~5/6
# padding vs. the stored <col width> attribute, so subtract it back off)
$wsCumplimiento.Columns("F").ColumnWidth = 24 - 5/6

$wsCumplimiento.Range("D3").Value = 8549.98
$wsCumplimiento.Range("E3").Value = 18907.0276
$wsCumplimiento.Range("F3").Value = 0.3113951864149974

$wsCumplimiento.Range("D12").Value = 414.97
$wsCumplimiento.Range("E12").Value = -64.97000000000003
$wsCumplimiento.Range("F12").Value = 1.185628571428571

$wsCumplimiento.Range("D19").Value = 46225.1
$wsCumplimiento.Range("E19").Value = 48222.34064517915
$wsCumplimiento.Range("F19").Value = 0.4894267084870918
